$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.249.53'
$ws.Range("E2").Value = '  -2.76%  '

$ws.Range("D3").Value = '3.555.59'
$ws.Range("E3").Value = '  -2.59%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '582.92'
$ws.Range("E5").Value = '  -0.98%  '

$ws.Range("D6").Value = '178.56'
$ws.Range("E6").Value = '  +0.12%  '

$ws.Range("D7").Value = '0.604'
$ws.Range("E7").Value = '  -3.69%  '

$ws.Range("E8").Value = '  +0.54%  '

$ws.Range("D9").Value = '0.664'
$ws.Range("E9").Value = '  -6.41%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  -11.21%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '53.14'
$ws.Range("E11").Value = '  -3.93%  '

$ws.Range("D12").Value = '0.0000249'
$ws.Range("E12").Value = '  -13.96%  '

$ws.Range("D13").Value = '9.81'
$ws.Range("E13").Value = '  -7.66%  '

$ws.Range("D14").Value = '4.143.90'
$ws.Range("E14").Value = '  -1.80%  '

$ws.Range("D15").Value = '3.558.40'
$ws.Range("E15").Value = '  -2.38%  '

$ws.Range("E16").Value = '  -0.42%  '

$ws.Range("D17").Value = '66.074.59'
$ws.Range("E17").Value = '  -2.65%  '

$ws.Range("D18").Value = '18.12'
$ws.Range("E18").Value = '  -5.82%  '

$ws.Range("D19").Value = '12.03'
$ws.Range("E19").Value = '  -5.07%  '

$ws.Range("E20").Value = '  -6.56%  '

$ws.Range("D21").Value = '388.08'
$ws.Range("E21").Value = '  -4.96%  '

$ws.Range("D22").Value = '4.24'
$ws.Range("E22").Value = '  -6.52%  '

$ws.Range("D23").Value = '83.80'
$ws.Range("E23").Value = '  -4.98%  '

$ws.Range("D24").Value = '2.83'
$ws.Range("E24").Value = '  -4.87%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '12.10'
$ws.Range("E25").Value = '  -3.81%  '

$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '6.00'
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("D27").Value = '10.11'
$ws.Range("E27").Value = '  -5.73%  '

$ws.Range("D28").Value = '3.60'
$ws.Range("E28").Value = '  -11.06%  '

$ws.Range("D29").Value = '8.84'
$ws.Range("E29").Value = '  -6.78%  '

$ws.Range("D30").Value = '30.80'
$ws.Range("E30").Value = '  -5.19%  '

$ws.Range("D31").Value = '6.67'
$ws.Range("E31").Value = '  -6.87%  '

$ws.Range("D32").Value = '65.10'
$ws.Range("E32").Value = '  +1.18%  '

$ws.Range("D33").Value = '11.78'
$ws.Range("E33").Value = '  -4.12%  '

$ws.Range("D34").Value = '588.92'
$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("D35").Value = '0.112'
$ws.Range("E35").Value = '  -4.64%  '

$ws.Range("D36").Value = '40.83'
$ws.Range("E36").Value = '  -4.09%  '

$ws.Range("E37").Value = '  -0.08%  '

$ws.Range("E38").Value = '  +0.07%  '

$ws.Range("D39").Value = '0.368'
$ws.Range("E39").Value = '  -6.70%  '

$ws.Range("D40").Value = '0.0₃0726'
$ws.Range("E40").Value = '  -17.11%  '

$ws.Range("D41").Value = '0.129'
$ws.Range("E41").Value = '  -5.50%  '

$ws.Range("D42").Value = '2.73'
$ws.Range("E42").Value = '  -8.87%  '

$ws.Range("D43").Value = '0.0405'
$ws.Range("E43").Value = '  -7.05%  '

$ws.Range("D44").Value = '2.736.90'
$ws.Range("E44").Value = '  +0.96%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '2.38'
$ws.Range("E45").Value = '  -11.22%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.08'
$ws.Range("E46").Value = '  -0.66%  '

$ws.Range("D47").Value = '0.129'
$ws.Range("E47").Value = '  -3.85%  '

$ws.Range("E48").Value = '  -7.11%  '

$ws.Range("D49").Value = '134.64'
$ws.Range("E49").Value = '  -3.57%  '

$ws.Range("E50").Value = '  -9.15%  '

$ws.Range("D51").Value = '2.54'
$ws.Range("E51").Value = '  -7.17%  '

